# Applies the "remoção de uma pagina" edit:
#  1. Five section headings change "Módulo" -> "Página" (keeping the
#     leading numbering and trailing qualifier text intact).
#  2. "Gerenciar perfis (Admin, Recepcionista, Podólogo)." drops "Podólogo".
#  3. The whole "8. Telas Adicionais (Opcionais)" sub-section (its leading
#     separator rule, heading and four body paragraphs) is removed.
#  4. The stray <w:lastRenderedPageBreak/> marker in front of
#     "Responsividade" is removed (content above it shrank, so the page
#     no longer breaks there) while the paragraph formatting is kept.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- 1 & 2: heading / body text tweaks -------------------------------
Replace-Text "3. Módulo de Agendamento" "3. Página de Agendamento"
Replace-Text "4. Módulo de Clientes" "4. Página de Clientes"
Replace-Text "4.3. Histórico do Cliente" "4.3. Página do Cliente"
Replace-Text "5. Módulo de Serviços e " "5. Página de Serviços e "
Replace-Text "6. Módulo de Relatórios" "6. Página de Relatórios"
Replace-Text "Gerenciar perfis (Admin, Recepcionista, Podólogo)." "Gerenciar perfis (Admin, Recepcionista)."

# --- 3: drop the "8. Telas Adicionais (Opcionais)" sub-section -------
# Locate the heading paragraph, then walk backwards to also take the
# horizontal-rule separator paragraph right before it, and forwards to
# take the four list paragraphs that belong to it. The following
# separator paragraph (the one introducing "Fluxo do Usuário") is left
# in place.
$count = $d.Paragraphs.Count
$headingIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "8. Telas Adicionais*") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -gt 0) {
    $startPara = $d.Paragraphs($headingIndex - 1)
    $endPara = $d.Paragraphs($headingIndex + 4)
    $blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $blockRange.Delete()
}

# --- 4: remove the stray lastRenderedPageBreak before "Responsividade"
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Responsividade*") {
        $full = $p.Range
        $r = $d.Range($full.Start, $full.End - 1)
        $runXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="404040"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>Responsividade</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $r.InsertXML($runXml)
        break
    }
}
